$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the two bullet-point paragraphs ("A 'randomly generated'
#    map ..." and "Turn system (play vs AI?) ...") plus the page-break
#    paragraph that used to separate them from the title page.
# ------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p3 = $d.Paragraphs.Item(3)
$introRange = $d.Range($p1.Range.Start, $p3.Range.End)
$introRange.Delete()

# ------------------------------------------------------------------
# 2. Retype the "FINAL PROJECT PROPOSAL" run so the stray
#    <w:lastRenderedPageBreak/> rendering artifact is dropped while
#    keeping the existing bold / language formatting.
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.Text = "FINAL PROJECT PROPOSAL"

# ------------------------------------------------------------------
# 3. Drop the leftover _GoBack bookmark just after the title.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 4. Append the new "Prototype program" section at the end of the
#    document.
# ------------------------------------------------------------------
$tail = $d.Content
$tail.Collapse(0)

$body = '<w:p><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>' + `
        '<w:p><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:rPr><w:b/><w:bCs/><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t>Prototype program</w:t></w:r></w:p>' + `
        '<w:p><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">The prototype for this project </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">focuses on generating a random map made of various tile types with weighted probabilities and </w:t></w:r></w:p>'

$xmlFrag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + "<w:body>$body</w:body>" + '</w:document></pkg:xmlData></pkg:part></pkg:package>'

$tail.InsertXML($xmlFrag)
